$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(312049950, Molham  Peretz: 9,3)"
$ws.Range("B1").Value = "(308073899, Anan  Kirshenbaum: 9,3)"
$ws.Range("C1").Value = "(318869187, Soaad  Leibovich: -9,-3)"
$ws.Range("D1").Value = "(205898513, Asaf  Braymok: -5,8)"
$ws.Range("E1").Value = "(316028364, Sami  Castro: -3,9)"
$ws.Range("F1").Value = "(318428158, Tal  Asulin: -2,9)"
$ws.Range("G1").Value = "(318294931, Shalev  Afanasenko: -1,1)"

$ws.Range("A3").Value = "cost: 328.0235462732961"
$ws.Range("A4").Value = "time: 60.60470925465924"
